$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------------------
# 1) Grow the table from 7 to 12 columns (still 4 rows for now, so the table
#    survives all the way through -- the engine drops a table whose range
#    would have zero data rows while resizing happens column-by-column).
# ---------------------------------------------------------------------------
$lo.Resize($ws.Range("A1:L4"))

# ---------------------------------------------------------------------------
# 2) Copy header formatting onto the cells that need it *before* we overwrite
#    the header text, so every header ends up using one of the two existing
#    "blue header" styles (no brand-new xf entries get minted).
#    A1 currently carries the numFmt "#,##0.00" blue header style; B1 carries
#    the numFmt "General" blue header style.
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("L1").PasteSpecial(-4122)

$ws.Range("B1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("A1").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Re-label the header row (new column layout).
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "FILIALE"
$ws.Range("B1").Value = "Réseau"
$ws.Range("C1").Value = "Type"
$ws.Range("D1").Value = "Date"
$ws.Range("E1").Value = "Devise"
$ws.Range("F1").Value = "Nbre Total De Transactions"
$ws.Range("G1").Value = "Montant Total de Transactions"
$ws.Range("H1").Value = "Montant de Transactions (Couverture)"
$ws.Range("I1").Value = "Nbre Total de Rejets"
$ws.Range("J1").Value = "Nbre de Transactions (Couverture)"
$ws.Range("K1").Value = "Rapprochement"
$ws.Range("L1").Value = "Montant de Rejets"

# ---------------------------------------------------------------------------
# 4) Drop the old sample data rows -- shrink the table to the header row only
#    first (so it stays alive with zero data rows), then delete the rows.
# ---------------------------------------------------------------------------
$lo.Resize($ws.Range("A1:L1"))
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()
$ws.Rows.Item(2).Delete()

# ---------------------------------------------------------------------------
# 5) New column widths.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 8.166666666666666
$ws.Columns.Item(2).ColumnWidth = 7.166666666666667
$ws.Columns.Item(3).ColumnWidth = 5.166666666666667
$ws.Columns.Item(4).ColumnWidth = 5.166666666666667
$ws.Columns.Item(5).ColumnWidth = 7.166666666666667
$ws.Columns.Item(6).ColumnWidth = 27.166666666666668
$ws.Columns.Item(7).ColumnWidth = 30.166666666666668
$ws.Columns.Item(8).ColumnWidth = 37.166666666666664
$ws.Columns.Item(9).ColumnWidth = 21.166666666666668
$ws.Columns.Item(10).ColumnWidth = 34.166666666666664
$ws.Columns.Item(11).ColumnWidth = 14.166666666666666
$ws.Columns.Item(12).ColumnWidth = 18.166666666666668
